# "chage type declation to follow julia rule"
#
# The header row on the "example1" sheet documented JSON array fields using
# a made-up "/path(Type)" notation. Update it to the Julia type-annotation
# style ("field::Vector{T}") used by the rest of the project.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("example1")

# Order matters here only insofar as it decides the order brand-new shared
# strings are appended in; write B1:D1 first, then A1, so the saved file's
# shared-string table ends with array_string, array_int, array_float, array_any
# (matching the canonical edit).
$ws1.Range("B1").Value = "/array_string::Vector{String}"
$ws1.Range("C1").Value = "/array_int::Vector{Int}"
$ws1.Range("D1").Value = "/array_float::Vector{Float64}"
$ws1.Range("A1").Value = "array_any::Vector"

# The author ended the editing session with "example1" on screen and cell
# D10 selected (instead of "example6" / C8 which was active beforehand).
$ws1.Activate() | Out-Null
$ws1.Range("D10").Select() | Out-Null
